# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This recomputes the "K" column (column G) values for each data row on
# Sheet1 and writes the refreshed results back into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed K values (s_vals), keyed by row number.
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 3
    6  = 1
    7  = 2
    8  = 1
    9  = 2
    10 = 0
    11 = 3
    12 = 0
    13 = 0
    14 = 2
    15 = 0
    16 = 2
    17 = 0
    18 = 2
    19 = 2
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 0
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Range("G$row").Value = $kValues[$row]
}
